# The author used "git jgit" to commit a change outside Webstudio: the
# greeting text in cell E8 ("Good Morning") was replaced with "GIT UPDATE".
# Apply that same edit here, then leave the active cell/selection on E8 to
# match the resulting workbook state (last-edited cell stays selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
